# Added clarity for Ketchup and Salsa.
#
# - Row 44 (die Salsa) -> rename the German word cell to "Salsa (Soßen)"
#   to disambiguate it from the condiment/app "Salsa".
# - Insert a new row right after the existing Ketchup row (old row 45)
#   for "der Ketchup" with an explicit masculine article, so both
#   article usages are represented; all following rows shift down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the Salsa entry for clarity.
$ws.Range("B44").Value = "Salsa (Soßen)"

# Insert a new row after row 45 (the existing "das Ketchup" row) and
# shift everything below down.
$ws.Rows.Item(46).Insert()

$ws.Range("A46").Value = "der"
$ws.Range("B46").Value = "Ketchup"
$ws.Range("C46").Value = "Ketchup"
$ws.Range("D46").Value = "Würzmittel und Soßen"
